$d = $word.ActiveDocument

# Locate the paragraph that holds the "{#tb5}{descripcion}{/tb5}" merge-field
# placeholder (directly under the "Objetivos Específicos" heading) using
# Find on a duplicated range so the original Content range is untouched.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("{#tb5}{descripcion}{/tb5}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the {#tb5}{descripcion}{/tb5} placeholder paragraph"
}

# Extend the range to swallow the trailing paragraph mark so the whole
# paragraph (not just its text) gets replaced by the new table.
$rng.MoveEnd(1, 1)

# Replace that paragraph with a one-column / two-row table: a blue header
# row reading "Objetivos Específicos" and a body row carrying the original
# {#tb5}...{/tb5} merge-field markup (split across runs, matching the
# target markup produced by Word's grammar-check proofing errors).
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblStyle w:val="Tablaconcuadrcula"/><w:tblW w:w="0" w:type="auto"/><w:jc w:val="center"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2658"/></w:tblGrid><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2658" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="365F91"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FFFFFF" w:themeColor="background1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Objetivos Específicos</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:jc w:val="center"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="2658" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>{#tb</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>5}</w:t></w:r><w:r><w:t>{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>descripcion</w:t></w:r><w:r><w:t>}{/tb</w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t>}</w:t></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$rng.InsertXML($xml)
